# Atmos 22 code now working
# Update Pin List sheet: append "/Teros21 ..." alt-sensor labels to the
# four 5TM connection descriptions (these sensors share the same pins).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pin List")

$ws.Range("C3").Value  = "Upper 5TM Data/Upper Teros21 Data"
$ws.Range("C5").Value  = "Upper 5TM Power/Upper Teros21 Power"
$ws.Range("C20").Value = "Lower 5TM Power/Lower Teros21 Power"
$ws.Range("C22").Value = "Lower 5TM Data/Lower Teros21 Data"

# Column C needs to widen to fit the new, longer labels (Excel auto-fit).
$ws.Columns.Item(3).AutoFit() | Out-Null
$ws.Columns.Item(3).ColumnWidth = 34.33

# Reset the view: scroll back to the top and move the active selection.
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("C23").Select()
